$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the min_val (D9) for the "reduction_factor" sensitivity row:
# widen the lower bound of the sensitivity sweep from 0.5 to 0.1.
$ws.Range("D9").Value = 0.1

# Move the active selection to match the author's last cursor position
# when the file was saved.
$null = $ws.Range("J10").Select()
